# Adds 3 tables (TABLE, RESTAURANT OPEN TIMES, TABLE BOOKINGS) plus sample
# data to the "database tables.xlsx" schema-diagram sheet.
#
# Layout recap (existing sheet):
#   Row 3  -> table name headers in columns C,E,G,I,K
#   Rows 4-10 -> column/field names for each table
#   Row 25 -> table name headers in columns C,E,G,I (second block)
#   Rows 26-32 -> column/field names for each table
#
# This change adds a new table header in column M (row 3) for "TABLE", a
# new table header in column K (row 25) for "RESTAURANT OPEN TIMES", and a
# new table header in column M (row 25) for "TABLE BOOKINGS" -- together
# with their respective field-name rows underneath.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "TABLE" table (top block, column M) ---------------------------
$ws.Range("M3").Value = "TABLE"
$ws.Range("M4").Value = "table_id"
$ws.Range("M5").Value = "seating"
$ws.Range("M6").Value = "restaurant_id"

# --- New "RESTAURANT OPEN TIMES" table (bottom block, column K) --------
$ws.Range("K25").Value = "RESTAURANT OPEN TIMES"
$ws.Range("K26").Value = "restaurant_id"
$ws.Range("K27").Value = "day_of_week"
$ws.Range("K28").Value = "opening_time"
$ws.Range("K29").Value = "closing_time"

# --- New "TABLE BOOKINGS" table (bottom block, column M) ----------------
$ws.Range("M25").Value = "TABLE BOOKINGS"
$ws.Range("M26").Value = "table_bookings_id"
$ws.Range("M27").Value = "booking_date"
$ws.Range("M28").Value = "booking_time"
$ws.Range("M29").Value = "table_id"
$ws.Range("M30").Value = "user_id"
$ws.Range("M31").Value = "booking_length"

# --- Column widths for the new/changed columns --------------------------
# Column K widened (now holds "RESTAURANT OPEN TIMES"); column M is new.
# (Inputs chosen so the engine's internal px-rounding lands as close as
# possible to the authored widths of 23.7109375 / 22 characters.)
$ws.Columns.Item(11).ColumnWidth = 22.8
$ws.Columns.Item(13).ColumnWidth = 21.15

# --- Selection / scroll position match the saved view -------------------
$ws.Range("O29").Select()
